# Courts.xlsx update:
#  - Add two new worksheets at the end: "CourtsApprovalAndAuditHistory" and
#    "CourtsAuditHistory", each with a header row + one data row mirroring
#    the style of the existing three tabs.
#  - Rename the TESTCASE value used across every tab from "test" to
#    "testT4149".
#  - Restore sheet selection/active-tab state: "Generate Document" becomes
#    the active tab (cell C18 selected) and "Petition Allegation" loses its
#    previously-active selection (D12 -> A2).

$wb = $excel.ActiveWorkbook

$wsGenerateDocument   = $wb.Worksheets.Item(1)
$wsGeneratePetition   = $wb.Worksheets.Item(2)
$wsPetitionAllegation = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------
# 1. Update the TESTCASE id (column A, row 2) on the three existing tabs.
# ---------------------------------------------------------------------
$wsGenerateDocument.Cells.Item(2, 1).Value   = "testT4149"
$wsGeneratePetition.Cells.Item(2, 1).Value   = "testT4149"
$wsPetitionAllegation.Cells.Item(2, 1).Value = "testT4149"

# ---------------------------------------------------------------------
# 2. Add the "CourtsApprovalAndAuditHistory" worksheet (after the last
#    existing sheet).
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsApproval = $wb.Worksheets.Add($null, $lastSheet)
$wsApproval.Name = "CourtsApprovalAndAuditHistory"

$approvalHeaders = @(
    "TESTCASE",
    "SCRIPT_ITERATION",
    "POM_ITERATION",
    "APPROVAL_AUDIT_HISTORY_TAB",
    "APPROVAL_HISTORY_TAB",
    "STEP_NAME_VERIFY",
    "STATUS_VERIFY",
    "ASSIGNED_TO_VERIFY",
    "ACTUAL_APPROVER_VERIFY",
    "COMMENTS_VERIFY",
    "COURT_WORK_ITEM_HISTORY_TAB",
    "DATE_VERIFY",
    "FIELD_VERIFY",
    "USER_VERIFY",
    "ORIGINAL_VALUE_VERIFY",
    "NEW_VALUE_VERIFY"
)
for ($c = 0; $c -lt $approvalHeaders.Length; $c++) {
    $wsApproval.Cells.Item(1, $c + 1).Value = $approvalHeaders[$c]
}

$approvalRow2 = @(
    "testT4149", 1, 1,
    "Click", "Click", "Yes", "Yes", "Yes", "Yes", "Yes",
    "Click", "Yes", "Yes", "Yes", "Yes", "Yes"
)
for ($c = 0; $c -lt $approvalRow2.Length; $c++) {
    $wsApproval.Cells.Item(2, $c + 1).Value = $approvalRow2[$c]
}

$wsApproval.Columns.Item(1).ColumnWidth  = 9.36
$wsApproval.Columns.Item(2).ColumnWidth  = 16.73
$wsApproval.Columns.Item(3).ColumnWidth  = 15.73
$wsApproval.Columns.Item(4).ColumnWidth  = 28.09
$wsApproval.Columns.Item(5).ColumnWidth  = 23.36
$wsApproval.Columns.Item(6).ColumnWidth  = 18.09
$wsApproval.Columns.Item(7).ColumnWidth  = 14.63
$wsApproval.Columns.Item(8).ColumnWidth  = 19.36
$wsApproval.Columns.Item(9).ColumnWidth  = 23.36
$wsApproval.Columns.Item(10).ColumnWidth = 18.18
$wsApproval.Columns.Item(11).ColumnWidth = 32
$wsApproval.Columns.Item(12).ColumnWidth = 14.27
$wsApproval.Columns.Item(13).ColumnWidth = 14.73
$wsApproval.Columns.Item(14).ColumnWidth = 14.91
$wsApproval.Columns.Item(15).ColumnWidth = 25.36

# Match the header styling convention used on the other tabs: column A's
# header label uses the explicit-black font, the rest stay default.
$wsApproval.Range("A1").Font.Color = 0

# ---------------------------------------------------------------------
# 3. Add the "CourtsAuditHistory" worksheet (after the approval-history
#    sheet, i.e. last tab overall).
# ---------------------------------------------------------------------
$wsAudit = $wb.Worksheets.Add($null, $wsApproval)
$wsAudit.Name = "CourtsAuditHistory"

$auditHeaders = @(
    "TESTCASE",
    "SCRIPT_ITERATION",
    "POM_ITERATION",
    "AUDIT_HISTORY_TAB",
    "COURT_WORK_ITEM_HISTORY_TAB",
    "DATE_VERIFY",
    "FIELD_VERIFY",
    "USER_VERIFY",
    "ORIGINAL_VALUE_VERIFY",
    "NEW_VALUE_VERIFY"
)
for ($c = 0; $c -lt $auditHeaders.Length; $c++) {
    $wsAudit.Cells.Item(1, $c + 1).Value = $auditHeaders[$c]
}

$auditRow2 = @(
    "testT4149", 1, 1,
    "Click", "Click", "Yes", "Yes", "Yes", "Yes", "Yes"
)
for ($c = 0; $c -lt $auditRow2.Length; $c++) {
    $wsAudit.Cells.Item(2, $c + 1).Value = $auditRow2[$c]
}

$wsAudit.Columns.Item(1).ColumnWidth  = 10.36
$wsAudit.Columns.Item(2).ColumnWidth  = 16.82
$wsAudit.Columns.Item(3).ColumnWidth  = 16.36
$wsAudit.Columns.Item(4).ColumnWidth  = 18.27
$wsAudit.Columns.Item(5).ColumnWidth  = 31.09
$wsAudit.Columns.Item(6).ColumnWidth  = 12.18
$wsAudit.Columns.Item(7).ColumnWidth  = 12.45
$wsAudit.Columns.Item(8).ColumnWidth  = 14.18
$wsAudit.Columns.Item(9).ColumnWidth  = 21.91
$wsAudit.Columns.Item(10).ColumnWidth = 17.63

# Same header styling convention as the approval-history tab.
$wsAudit.Range("A1").Font.Color = 0

# ---------------------------------------------------------------------
# 4. Restore the view/selection state seen in the final workbook:
#    - Petition Allegation: no longer the active tab, selection back to A2.
#    - CourtsAuditHistory: scrolled slightly, E16 selected.
#    - CourtsApprovalAndAuditHistory: default A2 selection.
#    - Generate Petition: default A2 selection (unchanged).
#    - Generate Document: becomes the active tab, C18 selected.
# ---------------------------------------------------------------------
$wsPetitionAllegation.Activate()
$wsPetitionAllegation.Range("A2").Select()

$wsAudit.Activate()
$wsAudit.Range("E16").Select()

$wsApproval.Activate()
$wsApproval.Range("A2").Select()

$wsGeneratePetition.Activate()
$wsGeneratePetition.Range("A2").Select()

$wsGenerateDocument.Activate()
$wsGenerateDocument.Range("C18").Select()
